# Merge the split text runs of the Title, Author and Abstract paragraphs
# into single runs containing the full paragraph text (no content change,
# just de-fragmenting the runs), per the commit's XML diff.

$d = $word.ActiveDocument

$targets = @{
    "Title"    = "Questions: Introduction to sigma notation"
    "Author"   = "Ifan Howells-Baines, Mark Toner"
    "Abstract" = "Questions relating to the guide on introduction to sigma notation."
}

foreach ($p in $d.Paragraphs) {
    $styleName = $p.Range.ParagraphStyle.NameLocal
    if ($targets.ContainsKey($styleName)) {
        $full = $targets[$styleName]
        $r = $p.Range
        $r.Find.Execute($full, $true, $false, $false, $false, $false,
                         $true, 1, $false, $full, 2)
    }
}
